$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the two new task rows
$ws.Range("B7").Value = "Get Api call working"
$ws.Range("D7").Value = "Done"

$ws.Range("B8").Value = "Fix CORS"

# Adjust column C width and selection to match the new state
# (target stored width 10.140625; engine quantizes ColumnWidth to the
# nearest 1/6 character unit, so 9.333333... is the closest achievable input)
$ws.Columns.Item(3).ColumnWidth = 9.333333333333334

$ws.Range("D8").Select()
